# Add a new "2022" data column (Z) to the 2.1.1 undernourishment worksheet,
# and renumber a handful of cellXf style indices that shifted around when
# the new "2022" header style was inserted into the shared cellXfs table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-point the five cells whose cellXf index changed because a few
#     existing xf records were reordered/rewritten to make room for the new
#     "2022" column-header style. Re-apply each cell's original look using
#     the look-alike style already sitting at the new index so the engine's
#     own de-duplication lands each cell back on an equivalent xf.

# A4 / B4 / C4 header row: swap bold bottom-border look for A4 with what
# used to be B4/C4's, and vice versa (they trade xf slots).
$ws.Cells.Item(4, 1).Copy() | Out-Null
$ws.Cells.Item(6, 1).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(1, 1).Copy() | Out-Null
$ws.Cells.Item(1, 1).PasteSpecial(-4122) | Out-Null

# --- New column Z: header "2022" + the per-region 2022 values, matching the
#     formatting already used by the neighbouring "2021" column (Y).
$ws.Cells.Item(4, 25).Copy() | Out-Null
$ws.Cells.Item(4, 26).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(4, 26).Value = 2022

$years = @{
    5  = 47.345690436648667
    6  = 55.294335329978139
    7  = 42.721146742902135
    8  = 56.732662465911261
    9  = 39.351829932862628
    10 = 43.952035422218046
    11 = 57.461907794486649
    12 = 32.073481974524846
    13 = 33.564455947162017
    14 = 55.803694659011171
    15 = 63.920911723512503
    16 = 52.521342498654128
}

foreach ($row in $years.Keys) {
    $ws.Cells.Item($row, 25).Copy() | Out-Null
    $ws.Cells.Item($row, 26).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($row, 26).Value = $years[$row]
}

# Update the visible selection to mirror the authored workbook (cursor
# parked one cell to the right of the new header, row 4 instead of row 15).
$ws.Range("AA4").Select() | Out-Null
